$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F:F").Select()
$ws.Range("F:F").Delete()
